# Update the daily cryptos snapshot (Price / Volume(1h) columns, plus the
# Maker/VeChain row swap at the bottom) to the new scrape.
#
# Price/Volume cells are stored as plain text (e.g. "67.401.49", "6.20",
# "  -2.61%  ") rather than numbers, so each literal is written with a
# leading apostrophe (Excel's standard "treat as text" quote-prefix) to
# stop COM from silently re-parsing them as numbers/percentages and
# mangling things like trailing zeros or thousand-dot groupings.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.401.49"
$ws.Range("E2").Value = "'  -2.61%  "
$ws.Range("D3").Value = "'3.719.12"
$ws.Range("E3").Value = "'  -3.26%  "
$ws.Range("E4").Value = "'  -0.12%  "
$ws.Range("D5").Value = "'598.32"
$ws.Range("E5").Value = "'  -0.79%  "
$ws.Range("D6").Value = "'166.99"
$ws.Range("E6").Value = "'  -3.46%  "
$ws.Range("D7").Value = "'3.713.92"
$ws.Range("E7").Value = "'  -3.33%  "
$ws.Range("E8").Value = "'  +0.03%  "
$ws.Range("D9").Value = "'0.534"
$ws.Range("E9").Value = "'  +1.35%  "
$ws.Range("D10").Value = "'0.166"
$ws.Range("E10").Value = "'  +1.58%  "
$ws.Range("D11").Value = "'6.20"
$ws.Range("E11").Value = "'  -2.81%  "
$ws.Range("E12").Value = "'  -3.30%  "
$ws.Range("D13").Value = "'37.85"
$ws.Range("E13").Value = "'  -3.14%  "
$ws.Range("D14").Value = "'0.0000243"
$ws.Range("E14").Value = "'  -2.64%  "
$ws.Range("D15").Value = "'4.341.42"
$ws.Range("E15").Value = "'  -3.21%  "
$ws.Range("D16").Value = "'3.714.73"
$ws.Range("E16").Value = "'  -3.27%  "
$ws.Range("D17").Value = "'67.419.97"
$ws.Range("E17").Value = "'  -2.88%  "
$ws.Range("D18").Value = "'7.28"
$ws.Range("E18").Value = "'  -1.83%  "
$ws.Range("D19").Value = "'17.50"
$ws.Range("E19").Value = "'  +7.08%  "
$ws.Range("E20").Value = "'  -2.21%  "
$ws.Range("D21").Value = "'487.72"
$ws.Range("E21").Value = "'  -2.55%  "
$ws.Range("D22").Value = "'9.27"
$ws.Range("E22").Value = "'  -2.92%  "
$ws.Range("D23").Value = "'0.729"
$ws.Range("E23").Value = "'  -2.06%  "
$ws.Range("D24").Value = "'85.09"
$ws.Range("E24").Value = "'  -2.86%  "
$ws.Range("D25").Value = "'0.0000144"
$ws.Range("E25").Value = "'  +4.01%  "
$ws.Range("E26").Value = "'  -4.59%  "
$ws.Range("D27").Value = "'12.23"
$ws.Range("E27").Value = "'  -2.39%  "
$ws.Range("D28").Value = "'10.05"
$ws.Range("E28").Value = "'  -2.19%  "
$ws.Range("E29").Value = "'  +0.11%  "
$ws.Range("E30").Value = "'  -1.73%  "
$ws.Range("D31").Value = "'2.36"
$ws.Range("E31").Value = "'  -5.54%  "
$ws.Range("D32").Value = "'7.70"
$ws.Range("E32").Value = "'  -1.36%  "
$ws.Range("D33").Value = "'31.52"
$ws.Range("E33").Value = "'  -4.72%  "
$ws.Range("D34").Value = "'3.859.23"
$ws.Range("E34").Value = "'  -3.24%  "
$ws.Range("E35").Value = "'  -3.86%  "
$ws.Range("D36").Value = "'3.662.15"
$ws.Range("E36").Value = "'  -2.95%  "
$ws.Range("D37").Value = "'0.998"
$ws.Range("E37").Value = "'  -0.11%  "
$ws.Range("D38").Value = "'0.996"
$ws.Range("E38").Value = "'  -3.54%  "
$ws.Range("D39").Value = "'5.83"
$ws.Range("E39").Value = "'  -3.07%  "
$ws.Range("D40").Value = "'0.132"
$ws.Range("E40").Value = "'  -5.29%  "
$ws.Range("E41").Value = "'  -2.31%  "
$ws.Range("D42").Value = "'48.76"
$ws.Range("E42").Value = "'  -1.50%  "
$ws.Range("D43").Value = "'428.59"
$ws.Range("E43").Value = "'  -7.35%  "
$ws.Range("E44").Value = "'  -1.48%  "
$ws.Range("E45").Value = "'  -5.13%  "
$ws.Range("E46").Value = "'  -0.40%  "
$ws.Range("E47").Value = "'  +0.02%  "
$ws.Range("D48").Value = "'40.58"
$ws.Range("E48").Value = "'  -3.82%  "
$ws.Range("D49").Value = "'140.57"
$ws.Range("E49").Value = "'  +1.14%  "
$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").Value = "'2.751.69"
$ws.Range("E50").Value = "'  -4.80%  "
$ws.Range("B51").Value = 'VeChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D51").Value = "'0.0351"
$ws.Range("E51").Value = "'  -2.27%  "
